$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D stays text (values look numeric but must remain exact strings, e.g. "35.314.50")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '35.314.50'

$ws.Range("D3").Value = '1.912.69'
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '0.724'
$ws.Range("E5").Value = '  +9.04%  '

$ws.Range("D6").Value = '255.58'
$ws.Range("E6").Value = '  +3.85%  '

$ws.Range("D8").Value = '40.76'
$ws.Range("E8").Value = '  -1.33%  '

$ws.Range("D9").Value = '0.372'
$ws.Range("E9").Value = '  +6.77%  '

$ws.Range("D10").Value = '52.80'
$ws.Range("E10").Value = '  +0.03%  '

$ws.Range("D11").Value = '0.0765'
$ws.Range("E11").Value = '  +6.74%  '

$ws.Range("D12").Value = '0.0988'

$ws.Range("D13").Value = '2.192.76'
$ws.Range("E13").Value = '  +0.43%  '

$ws.Range("D14").Value = '12.85'
$ws.Range("E14").Value = '  +6.22%  '

$ws.Range("E15").Value = '  +4.14%  '

$ws.Range("D16").Value = '4.97'
$ws.Range("E16").Value = '  +2.33%  '

$ws.Range("D17").Value = '1.935.18'
$ws.Range("E17").Value = '  +2.16%  '

$ws.Range("D18").Value = '35.317.35'
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("D19").Value = '74.59'
$ws.Range("E19").Value = '  +2.93%  '

$ws.Range("D20").Value = '0.0₃0852'
$ws.Range("E20").Value = '  +2.93%  '

$ws.Range("D21").Value = '244.37'
$ws.Range("E21").Value = '  +1.99%  '

$ws.Range("D22").Value = '13.10'
$ws.Range("E22").Value = '  +4.92%  '

$ws.Range("D23").Value = '5.12'
$ws.Range("E23").Value = '  +5.75%  '

$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  +4.46%  '

$ws.Range("D26").Value = '2.43'
$ws.Range("E26").Value = '  +3.09%  '

$ws.Range("D27").Value = '166.57'
$ws.Range("E27").Value = '  -2.12%  '

$ws.Range("D28").Value = '8.70'
$ws.Range("E28").Value = '  +3.08%  '

$ws.Range("D29").Value = '18.76'
$ws.Range("E29").Value = '  +1.92%  '

$ws.Range("E30").Value = '  +4.59%  '

$ws.Range("D31").Value = '4.128.92'
$ws.Range("E31").Value = '  +19.46%  '

$ws.Range("E32").Value = '  +5.42%  '

$ws.Range("E33").Value = '  +13.95%  '

$ws.Range("E34").Value = '  +24.00%  '

$ws.Range("E35").Value = '  +4.20%  '

$ws.Range("E36").Value = '  +3.69%  '

$ws.Range("E37").Value = '  -0.95%  '

$ws.Range("D38").Value = '0.913'
$ws.Range("E38").Value = '  -2.14%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").Value = '17.27'
$ws.Range("E40").Value = '  +5.58%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0219'
$ws.Range("E41").Value = '  +5.11%  '

$ws.Range("D42").Value = '97.01'
$ws.Range("E42").Value = '  +7.85%  '

$ws.Range("E43").Value = '  +1.41%  '

$ws.Range("D44").Value = '0.0651'
$ws.Range("E44").Value = '  +2.30%  '

$ws.Range("D45").Value = '1.337.83'

$ws.Range("E46").Value = '  +2.65%  '

$ws.Range("E47").Value = '  +0.83%  '

$ws.Range("E48").Value = '  +3.20%  '

$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").Value = '45.23'
$ws.Range("E50").Value = '  -5.06%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0748'
$ws.Range("E51").Value = '  +5.78%  '
